$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, date range) ---
$ws.Range("A8").Value = "Volume 32   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/17/2025  Through  11/23/2025"

# --- Fix number formats for cells that flip between text placeholder and numeric ---
# Reference cells with the desired target style already applied:
#   style "14" (integer #,##0)      -> C15
#   style "15" (percent #,##0.0)    -> E15
#   style "13" (text placeholder)   -> C14
$ws.Range("C15").Copy() | Out-Null
$ws.Range("D14,D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").Copy() | Out-Null
$ws.Range("E14,E22").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("F14,D33,E33").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Data value updates (rows 14-33, 40-46) ---
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = -100
$ws.Range("F14").Value = "0"
$ws.Range("H14").Value = -100
$ws.Range("J14").Value = 7
$ws.Range("K14").Value = 71.428571428571
$ws.Range("F15").Value = 7
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = 40
$ws.Range("I15").Value = 56
$ws.Range("J15").Value = 40
$ws.Range("K15").Value = 40
$ws.Range("L15").Value = 51.351351351351
$ws.Range("M15").Value = 51.351351351351
$ws.Range("N15").Value = -17.647058823529
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = -45.454545454545
$ws.Range("F16").Value = 38
$ws.Range("G16").Value = 46
$ws.Range("H16").Value = -17.391304347826
$ws.Range("I16").Value = 442
$ws.Range("J16").Value = 489
$ws.Range("K16").Value = -9.61145194274
$ws.Range("L16").Value = 3.7558685446
$ws.Range("M16").Value = 11.055276381909
$ws.Range("N16").Value = -64.64
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = -22.222222222222
$ws.Range("G17").Value = 69
$ws.Range("H17").Value = -17.391304347826
$ws.Range("I17").Value = 908
$ws.Range("J17").Value = 738
$ws.Range("K17").Value = 23.035230352303
$ws.Range("L17").Value = 15.668789808917
$ws.Range("M17").Value = 130.456852791878
$ws.Range("N17").Value = 16.410256410256
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -37.5
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -22.727272727272
$ws.Range("I18").Value = 231
$ws.Range("J18").Value = 241
$ws.Range("K18").Value = -4.149377593361
$ws.Range("L18").Value = -1.702127659574
$ws.Range("M18").Value = -30.211480362537
$ws.Range("N18").Value = -86.323268206039
$ws.Range("C19").Value = 23
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 76.923076923076
$ws.Range("F19").Value = 71
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = 24.561403508771
$ws.Range("I19").Value = 909
$ws.Range("J19").Value = 799
$ws.Range("K19").Value = 13.767209011264
$ws.Range("L19").Value = 28.571428571428
$ws.Range("M19").Value = 225.806451612903
$ws.Range("N19").Value = 82.897384305835
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 14
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 34
$ws.Range("G20").Value = 41
$ws.Range("H20").Value = -17.073170731707
$ws.Range("I20").Value = 569
$ws.Range("J20").Value = 499
$ws.Range("K20").Value = 14.028056112224
$ws.Range("L20").Value = -1.043478260869
$ws.Range("M20").Value = 85.9477124183
$ws.Range("N20").Value = -61.027397260274
$ws.Range("C21").Value = 56
$ws.Range("D21").Value = 66
$ws.Range("E21").Value = -15.151515151515
$ws.Range("F21").Value = 224
$ws.Range("G21").Value = 241
$ws.Range("H21").Value = -7.053941908713
$ws.Range("I21").Value = 3127
$ws.Range("J21").Value = 2813
$ws.Range("K21").Value = 11.162460007109
$ws.Range("L21").Value = 12.684684684684
$ws.Range("M21").Value = 76.966610073571
$ws.Range("N21").Value = -45.82467082467
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -60
$ws.Range("I22").Value = 23
$ws.Range("J22").Value = 28
$ws.Range("K22").Value = -17.857142857142
$ws.Range("L22").Value = -11.538461538461
$ws.Range("M22").Value = -8
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 400
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 85.714285714285
$ws.Range("I23").Value = 150
$ws.Range("J23").Value = 109
$ws.Range("K23").Value = 37.614678899082
$ws.Range("L23").Value = 61.290322580645
$ws.Range("M23").Value = 120.588235294118
$ws.Range("C24").Value = 34
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = 6.25
$ws.Range("F24").Value = 108
$ws.Range("G24").Value = 136
$ws.Range("H24").Value = -20.588235294117
$ws.Range("I24").Value = 1327
$ws.Range("J24").Value = 1193
$ws.Range("K24").Value = 11.232187761944
$ws.Range("L24").Value = 11.606391925988
$ws.Range("M24").Value = 92.597968069666
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 266.666666666667
$ws.Range("F25").Value = 25
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = -28.571428571428
$ws.Range("I25").Value = 201
$ws.Range("J25").Value = 352
$ws.Range("K25").Value = -42.897727272727
$ws.Range("L25").Value = -37.770897832817
$ws.Range("C26").Value = 27
$ws.Range("D26").Value = 21
$ws.Range("E26").Value = 28.571428571428
$ws.Range("F26").Value = 87
$ws.Range("G26").Value = 78
$ws.Range("H26").Value = 11.538461538461
$ws.Range("I26").Value = 1099
$ws.Range("J26").Value = 1009
$ws.Range("K26").Value = 8.919722497522
$ws.Range("L26").Value = 23.206278026905
$ws.Range("M26").Value = 25.028441410694
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -22.222222222222
$ws.Range("I27").Value = 72
$ws.Range("J27").Value = 60
$ws.Range("K27").Value = 20
$ws.Range("L27").Value = 18.032786885245
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 600
$ws.Range("I28").Value = 99
$ws.Range("K28").Value = 32
$ws.Range("L28").Value = 37.5
$ws.Range("G29").Value = 3
$ws.Range("J29").Value = 32
$ws.Range("K29").Value = 18.75
$ws.Range("M29").Value = -39.682539682539
$ws.Range("N29").Value = -69.10569105691
$ws.Range("G30").Value = 3
$ws.Range("J30").Value = 23
$ws.Range("K30").Value = 17.391304347826
$ws.Range("M30").Value = -47.058823529411
$ws.Range("N30").Value = -76.724137931034
$ws.Range("D33").Value = "0"
$ws.Range("E33").Value = "***.*"
$ws.Range("L33").Value = 0
